$wb = $excel.ActiveWorkbook

# --- Sheet "Dry Cakes": add 5 new Cheesecake rows (10-14) ---
$dry = $wb.Worksheets.Item("Dry Cakes")

# Replicate the existing row-9 (dc8) formatting pattern into the new rows
# so that columns A, B, G, H keep style index 2 while C, D, E, F stay
# unstyled, matching the layout already used by rows 2-9.
$dry.Range("A9:H9").Copy()
$dry.Range("A10:H10").PasteSpecial(-4122)
$dry.Range("A11:H11").PasteSpecial(-4122)
$dry.Range("A12:H12").PasteSpecial(-4122)
$dry.Range("A13:H13").PasteSpecial(-4122)
$dry.Range("A14:H14").PasteSpecial(-4122)
$dry.Range("D10:E14").Clear()

# Names (column B)
$dry.Range("B10").Value = "New York Cheesecake"
$dry.Range("B11").Value = "Lemon  Cheesecake"
$dry.Range("B12").Value = "Strawberry  Cheesecake"
$dry.Range("B13").Value = "Blueberry  Cheesecake"
$dry.Range("B14").Value = "Biscoff  Cheesecake"

# Ids (column A)
$dry.Range("A10").Value = "dc9"
$dry.Range("A11").Value = "dc10"
$dry.Range("A12").Value = "dc11"
$dry.Range("A13").Value = "dc12"
$dry.Range("A14").Value = "dc13"

# Images (column C)
$dry.Range("C11").Value = "dry-cakes/Lemon-Cheesecake.jpg"
$dry.Range("C12").Value = "dry-cakes/Strawberry-Cheesecake.jpg"
$dry.Range("C13").Value = "dry-cakes/Blueberry-Cheesecake.jpg"
$dry.Range("C14").Value = "dry-cakes/Biscoff-Cheesecake.jpg"
$dry.Range("C10").Value = "dry-cakes/New-York-Cheesecake.jpg"

# Description (column F)
$dry.Range("F10").Value = "[Veg preparation]"
$dry.Range("F11").Value = "[Veg preparation]"
$dry.Range("F12").Value = "[Veg preparation]"
$dry.Range("F13").Value = "[Veg preparation]"
$dry.Range("F14").Value = "[Veg preparation]"

# inStock (column G)
$dry.Range("G10").Value = "yes"
$dry.Range("G11").Value = "yes"
$dry.Range("G12").Value = "yes"
$dry.Range("G13").Value = "yes"
$dry.Range("G14").Value = "yes"

# onDiscount (column H)
$dry.Range("H10").Value = "no"
$dry.Range("H11").Value = "no"
$dry.Range("H12").Value = "no"
$dry.Range("H13").Value = "no"
$dry.Range("H14").Value = "no"

# --- Switch the active sheet/selection from "Sweet Delights" to "Dry Cakes" ---
# (Sweet Delights keeps its previously cached D5 selection; it just stops
# being the active/selected tab once Dry Cakes is activated below.)
$dry.Activate()
$dry.Range("H15").Select()
